$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new (blank) rows at position 7. This pushes the existing
#     rows 7/8 (the 2021-02-11 / 44238 week) down to rows 9/10, preserving
#     their original values untouched. ---
$ws.Rows.Item(7).Resize(2).Insert()

# --- Row 5: Primera, week updated from 2021-02-15 (44242) to a new week (44424) ---
$ws.Cells.Item(5, 4).Value = 44424        # D5 Fecha
$ws.Cells.Item(5, 10).Value = 75          # J5 Volumen
$ws.Cells.Item(5, 11).Value = 18000       # K5 Precio minimo
$ws.Cells.Item(5, 12).Value = 18000       # L5 Precio maximo
$ws.Cells.Item(5, 13).Value = 18000       # M5 Precio promedio ponderado
$ws.Cells.Item(5, 14).Value = "$/caja 15 kilos"   # N5 Unidad de comercializacion
$ws.Cells.Item(5, 16).Value = 1200        # P5 Precio $/Kg
$ws.Cells.Item(5, 17).Value = 15          # Q5 Kg o Unidades

# --- Row 6: Segunda, week updated from 2021-02-15 (44242) to a new week (44424) ---
$ws.Cells.Item(6, 4).Value = 44424        # D6 Fecha
$ws.Cells.Item(6, 11).Value = 12000       # K6 Precio minimo
$ws.Cells.Item(6, 12).Value = 12000       # L6 Precio maximo
$ws.Cells.Item(6, 13).Value = 12000       # M6 Precio promedio ponderado
$ws.Cells.Item(6, 14).Value = "$/caja 15 kilos"   # N6 Unidad de comercializacion
$ws.Cells.Item(6, 16).Value = 800         # P6 Precio $/Kg
$ws.Cells.Item(6, 17).Value = 15          # Q6 Kg o Unidades

# --- Row 7 (newly inserted, blank): Primera, week 2021-02-15 (44242) ---
$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44242
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 100112043
$ws.Cells.Item(7, 7).Value = "Pepino dulce"
$ws.Cells.Item(7, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 60
$ws.Cells.Item(7, 11).Value = 13000
$ws.Cells.Item(7, 12).Value = 13000
$ws.Cells.Item(7, 13).Value = 13000
$ws.Cells.Item(7, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 722
$ws.Cells.Item(7, 17).Value = 18
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# --- Row 8 (newly inserted, blank): Segunda, week 2021-02-15 (44242) ---
$ws.Cells.Item(8, 1).Value = 3
$ws.Cells.Item(8, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 44242
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 100112043
$ws.Cells.Item(8, 7).Value = "Pepino dulce"
$ws.Cells.Item(8, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(8, 9).Value = "Segunda"
$ws.Cells.Item(8, 10).Value = 50
$ws.Cells.Item(8, 11).Value = 10000
$ws.Cells.Item(8, 12).Value = 10000
$ws.Cells.Item(8, 13).Value = 10000
$ws.Cells.Item(8, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 556
$ws.Cells.Item(8, 17).Value = 18
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# Rows 9/10 already hold the original (untouched) 2021-02-11 data because the
# Insert() above shifted the former rows 7/8 down automatically.

# Make sure the Fecha (date) cells for the newly-inserted rows 7/8 use the
# same date/time number format as the rest of column D.
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
